$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = '''2019-03-12'
$ws.Cells.Item(3, 6).Value = '''2019-03-12'
$ws.Cells.Item(4, 5).Value = 16
$ws.Cells.Item(4, 6).Value = '''2019-03-12'
$ws.Cells.Item(5, 5).Value = 52
$ws.Cells.Item(5, 6).Value = '''2019-03-12'
$ws.Cells.Item(6, 6).Value = '''2019-03-12'
$ws.Cells.Item(7, 6).Value = '''2019-03-12'
$ws.Cells.Item(8, 6).Value = '''2019-03-12'
$ws.Cells.Item(9, 3).Value = 'LOOSE CHANGE'
$ws.Cells.Item(9, 4).Value = 1507
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = '''2019-03-12'
$ws.Cells.Item(10, 3).Value = 'HOLIDAY LOOSE CHANGE'
$ws.Cells.Item(10, 4).Value = 1514
$ws.Cells.Item(10, 5).Value = 5
$ws.Cells.Item(10, 6).Value = '''2019-03-12'
$ws.Cells.Item(11, 6).Value = '''2019-03-12'
$ws.Cells.Item(12, 6).Value = '''2019-03-12'
$ws.Cells.Item(13, 6).Value = '''2019-03-12'
$ws.Cells.Item(14, 6).Value = '''2019-03-12'
$ws.Cells.Item(15, 6).Value = '''2019-03-12'
$ws.Cells.Item(16, 6).Value = '''2019-03-12'
$ws.Cells.Item(17, 6).Value = '''2019-03-12'
$ws.Cells.Item(18, 6).Value = '''2019-03-12'
$ws.Cells.Item(19, 6).Value = '''2019-03-12'
$ws.Cells.Item(20, 6).Value = '''2019-03-12'
$ws.Cells.Item(21, 6).Value = '''2019-03-12'
$ws.Cells.Item(22, 6).Value = '''2019-03-12'
$ws.Cells.Item(23, 6).Value = '''2019-03-12'
$ws.Cells.Item(24, 3).Value = 'FAST CASH'
$ws.Cells.Item(24, 4).Value = 1518
$ws.Cells.Item(24, 5).Value = 1
$ws.Cells.Item(24, 6).Value = '''2019-03-12'
$ws.Cells.Item(25, 3).Value = 'BIRTHDAY SURPRISE'
$ws.Cells.Item(25, 4).Value = 1524
$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 6).Value = '''2019-03-12'
$ws.Cells.Item(26, 3).Value = 'MONEY ROLL'
$ws.Cells.Item(26, 4).Value = 1541
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = '''2019-03-12'
$ws.Cells.Item(27, 3).Value = 'FAST MONEY'
$ws.Cells.Item(27, 4).Value = 1548
$ws.Cells.Item(27, 5).Value = 3
$ws.Cells.Item(27, 6).Value = '''2019-03-12'
$ws.Cells.Item(28, 6).Value = '''2019-03-12'
$ws.Cells.Item(29, 6).Value = '''2019-03-12'
$ws.Cells.Item(30, 6).Value = '''2019-03-12'
$ws.Cells.Item(31, 3).Value = 'Holiday Lucky Times 10'
$ws.Cells.Item(31, 4).Value = 1450
$ws.Cells.Item(31, 6).Value = '''2019-03-12'
$ws.Cells.Item(32, 3).Value = 'BIG MONEY SPECTACULAR'
$ws.Cells.Item(32, 4).Value = 1486
$ws.Cells.Item(32, 6).Value = '''2019-03-12'
$ws.Cells.Item(33, 3).Value = '$20,000 GOLD RUSH'
$ws.Cells.Item(33, 4).Value = 1468
$ws.Cells.Item(33, 6).Value = '''2019-03-12'
$ws.Cells.Item(34, 3).Value = 'SUMMER LUCKY TIMES 10'
$ws.Cells.Item(34, 4).Value = 1487
$ws.Cells.Item(34, 6).Value = '''2019-03-12'
$ws.Cells.Item(35, 3).Value = '10X THE BUCKS'
$ws.Cells.Item(35, 4).Value = 1442
$ws.Cells.Item(35, 6).Value = '''2019-03-12'
$ws.Cells.Item(36, 3).Value = 'HIT $100!'
$ws.Cells.Item(36, 4).Value = 1445
$ws.Cells.Item(36, 6).Value = '''2019-03-12'
$ws.Cells.Item(37, 3).Value = 'MONEY ROLL'
$ws.Cells.Item(37, 4).Value = 1500
$ws.Cells.Item(37, 6).Value = '''2019-03-12'
$ws.Cells.Item(38, 6).Value = '''2019-03-12'
$ws.Cells.Item(39, 6).Value = '''2019-03-12'
$ws.Cells.Item(40, 6).Value = '''2019-03-12'
$ws.Cells.Item(41, 6).Value = '''2019-03-12'
$ws.Cells.Item(42, 6).Value = '''2019-03-12'
$ws.Cells.Item(43, 6).Value = '''2019-03-12'
$ws.Cells.Item(44, 6).Value = '''2019-03-12'
$ws.Cells.Item(45, 6).Value = '''2019-03-12'
$ws.Cells.Item(46, 6).Value = '''2019-03-12'
$ws.Cells.Item(47, 6).Value = '''2019-03-12'
$ws.Cells.Item(48, 3).Value = 'BINGO PLUS'
$ws.Cells.Item(48, 4).Value = 1508
$ws.Cells.Item(48, 5).Value = 4
$ws.Cells.Item(48, 6).Value = '''2019-03-12'
$ws.Cells.Item(49, 3).Value = 'CROSSWORD'
$ws.Cells.Item(49, 4).Value = 1438
$ws.Cells.Item(49, 5).Value = 1
$ws.Cells.Item(49, 6).Value = '''2019-03-12'
$ws.Cells.Item(50, 6).Value = '''2019-03-12'
$ws.Cells.Item(51, 3).Value = 'LOTERIA'
$ws.Cells.Item(51, 4).Value = 1409
$ws.Cells.Item(51, 6).Value = '''2019-03-12'
$ws.Cells.Item(52, 3).Value = 'CROSSWORD'
$ws.Cells.Item(52, 4).Value = 1448
$ws.Cells.Item(52, 6).Value = '''2019-03-12'
$ws.Cells.Item(53, 6).Value = '''2019-03-12'
$ws.Cells.Item(54, 6).Value = '''2019-03-12'
$ws.Cells.Item(55, 6).Value = '''2019-03-12'
$ws.Cells.Item(56, 3).Value = '$500 FRENZY'
$ws.Cells.Item(56, 4).Value = 1512
$ws.Cells.Item(56, 5).Value = 2
$ws.Cells.Item(56, 6).Value = '''2019-03-12'
$ws.Cells.Item(57, 3).Value = 7
$ws.Cells.Item(57, 4).Value = 1457
$ws.Cells.Item(57, 6).Value = '''2019-03-12'
$ws.Cells.Item(58, 3).Value = 'Holiday Spectacular'
$ws.Cells.Item(58, 4).Value = 1451
$ws.Cells.Item(58, 5).Value = 1
$ws.Cells.Item(58, 6).Value = '''2019-03-12'
$ws.Cells.Item(59, 6).Value = '''2019-03-12'
$ws.Cells.Item(60, 3).Value = 'SUPER CROSSWORD'
$ws.Cells.Item(60, 4).Value = 1526
$ws.Cells.Item(60, 5).Value = 6
$ws.Cells.Item(60, 6).Value = '''2019-03-12'
$ws.Cells.Item(61, 3).Value = 'HOLIDAY BUCKS'
$ws.Cells.Item(61, 4).Value = 1516
$ws.Cells.Item(61, 5).Value = 1
$ws.Cells.Item(61, 6).Value = '''2019-03-12'
$ws.Cells.Item(62, 6).Value = '''2019-03-12'
$ws.Cells.Item(63, 6).Value = '''2019-03-12'
$ws.Cells.Item(64, 3).Value = '$200,000 LUCKY 7S'
$ws.Cells.Item(64, 4).Value = 1504
$ws.Cells.Item(64, 6).Value = '''2019-03-12'
$ws.Cells.Item(65, 3).Value = 'Diamonds and Dollars'
$ws.Cells.Item(65, 4).Value = 1446
$ws.Cells.Item(65, 6).Value = '''2019-03-12'
$ws.Cells.Item(66, 6).Value = '''2019-03-12'
$ws.Cells.Item(67, 6).Value = '''2019-03-12'
$ws.Cells.Item(68, 3).Value = '$150K POKER RICHE$'
$ws.Cells.Item(68, 4).Value = 1523
$ws.Cells.Item(68, 5).Value = 3
$ws.Cells.Item(68, 6).Value = '''2019-03-12'
$ws.Cells.Item(69, 3).Value = '$200,000 GOLD RUSH'
$ws.Cells.Item(69, 4).Value = 1469
$ws.Cells.Item(69, 5).Value = 1
$ws.Cells.Item(69, 6).Value = '''2019-03-12'
$ws.Cells.Item(70, 6).Value = '''2019-03-12'
$ws.Cells.Item(71, 3).Value = '$18,000,000 SILVER PAYOUT'
$ws.Cells.Item(71, 4).Value = 1531
$ws.Cells.Item(71, 5).Value = 3
$ws.Cells.Item(71, 6).Value = '''2019-03-12'
$ws.Cells.Item(72, 3).Value = 'LADY LUCK BONUS'
$ws.Cells.Item(72, 4).Value = 1497
$ws.Cells.Item(72, 5).Value = 2
$ws.Cells.Item(72, 6).Value = '''2019-03-12'
$ws.Cells.Item(73, 6).Value = '''2019-03-12'
$ws.Cells.Item(74, 6).Value = '''2019-03-12'
$ws.Cells.Item(75, 6).Value = '''2019-03-12'
$ws.Cells.Item(76, 6).Value = '''2019-03-12'
$ws.Cells.Item(77, 6).Value = '''2019-03-12'
$ws.Cells.Item(78, 3).Value = '$150,000 POKER'
$ws.Cells.Item(78, 4).Value = 1455
$ws.Cells.Item(78, 6).Value = '''2019-03-12'
$ws.Cells.Item(79, 3).Value = '$500 FRENZY'
$ws.Cells.Item(79, 4).Value = 1466
$ws.Cells.Item(79, 6).Value = '''2019-03-12'
$ws.Cells.Item(80, 3).Value = 'Hot Spot Bingo'
$ws.Cells.Item(80, 4).Value = 1414
$ws.Cells.Item(80, 6).Value = '''2019-03-12'
$ws.Cells.Item(81, 3).Value = 'SUPER CROSSWORD'
$ws.Cells.Item(81, 4).Value = 1472
$ws.Cells.Item(81, 6).Value = '''2019-03-12'
$ws.Cells.Item(82, 3).Value = 'Bingo Times 10'
$ws.Cells.Item(82, 4).Value = 1501
$ws.Cells.Item(82, 6).Value = '''2019-03-12'
$ws.Cells.Item(83, 3).Value = 'Super Crossword'
$ws.Cells.Item(83, 4).Value = 1477
$ws.Cells.Item(83, 6).Value = '''2019-03-12'
$ws.Cells.Item(84, 3).Value = 'SUPER CROSSWORD'
$ws.Cells.Item(84, 4).Value = 1509
$ws.Cells.Item(84, 6).Value = '''2019-03-12'
$ws.Cells.Item(85, 3).Value = '20X THE BUCKS'
$ws.Cells.Item(85, 4).Value = 1443
$ws.Cells.Item(85, 6).Value = '''2019-03-12'
$ws.Cells.Item(86, 6).Value = '''2019-03-12'
$ws.Cells.Item(87, 6).Value = '''2019-03-12'
$ws.Cells.Item(88, 3).Value = '50X THE BUCKS'
$ws.Cells.Item(88, 4).Value = 1444
$ws.Cells.Item(88, 6).Value = '''2019-03-12'
$ws.Cells.Item(89, 3).Value = '$1,000,000 SPECTACULAR'
$ws.Cells.Item(89, 4).Value = 1463
$ws.Cells.Item(89, 6).Value = '''2019-03-12'
$ws.Cells.Item(90, 3).Value = '$1,000,000 LUCKY 7S'
$ws.Cells.Item(90, 4).Value = 1505
$ws.Cells.Item(90, 5).Value = 3
$ws.Cells.Item(90, 6).Value = '''2019-03-12'
$ws.Cells.Item(91, 3).Value = 'HOT RICHES'
$ws.Cells.Item(91, 4).Value = 1489
$ws.Cells.Item(91, 5).Value = 1
$ws.Cells.Item(91, 6).Value = '''2019-03-12'
$ws.Cells.Item(92, 3).Value = '$1,000,000 Spectacular'
$ws.Cells.Item(92, 4).Value = 1368
$ws.Cells.Item(93, 6).Value = '''2019-03-12'
$ws.Cells.Item(94, 6).Value = '''2019-03-12'
$ws.Cells.Item(95, 6).Value = '''2019-03-12'
$ws.Cells.Item(96, 6).Value = '''2019-03-12'
$ws.Cells.Item(97, 3).Value = '$250,000 CROSSWORD'
$ws.Cells.Item(97, 4).Value = 1542
$ws.Cells.Item(97, 5).Value = 2
$ws.Cells.Item(97, 6).Value = '''2019-03-12'
$ws.Cells.Item(98, 3).Value = 'LUCKY TIMES 50'
$ws.Cells.Item(98, 4).Value = 1536
$ws.Cells.Item(98, 5).Value = 3
$ws.Cells.Item(98, 6).Value = '''2019-03-12'
$ws.Cells.Item(99, 5).Value = 462
$ws.Cells.Item(99, 6).Value = '''2019-03-12'
$ws.Cells.Item(100, 3).Value = '$1,000,000 GOLD RUSH'
$ws.Cells.Item(100, 4).Value = 1470
$ws.Cells.Item(100, 5).Value = 1
$ws.Cells.Item(100, 6).Value = '''2019-03-12'
$ws.Cells.Item(101, 3).Value = 'JERSEY CASH BLOWOUT'
$ws.Cells.Item(101, 4).Value = 1482
$ws.Cells.Item(101, 5).Value = 31231
$ws.Cells.Item(101, 6).Value = '''2019-03-12'
$ws.Cells.Item(102, 6).Value = '''2019-03-12'
$ws.Cells.Item(103, 6).Value = '''2019-03-12'
$ws.Cells.Item(104, 4).Value = 1528
$ws.Cells.Item(104, 6).Value = '''2019-03-12'
$ws.Cells.Item(105, 4).Value = 1458
$ws.Cells.Item(105, 6).Value = '''2019-03-12'
$ws.Cells.Item(106, 4).Value = 1473
$ws.Cells.Item(106, 6).Value = '''2019-03-12'
$ws.Cells.Item(107, 4).Value = 1478
$ws.Cells.Item(107, 6).Value = '''2019-03-12'
$ws.Cells.Item(108, 4).Value = 1493
$ws.Cells.Item(108, 6).Value = '''2019-03-12'
$ws.Cells.Item(109, 6).Value = '''2019-03-12'
$ws.Cells.Item(110, 6).Value = '''2019-03-12'
$ws.Cells.Item(111, 6).Value = '''2019-03-12'
$ws.Cells.Item(112, 6).Value = '''2019-03-12'
$ws.Cells.Item(113, 6).Value = '''2019-03-12'
$ws.Cells.Item(114, 6).Value = '''2019-03-12'
$ws.Cells.Item(115, 3).Value = 'Platinum Diamond Spectacular'
$ws.Cells.Item(115, 4).Value = 1364
$ws.Cells.Item(115, 6).Value = '''2019-03-12'
$ws.Cells.Item(116, 3).Value = 'MEGA CROSSWORD'
$ws.Cells.Item(116, 4).Value = 1413
$ws.Cells.Item(116, 6).Value = '''2019-03-12'
$ws.Cells.Item(117, 3).Value = '100X THE BUCKS'
$ws.Cells.Item(117, 4).Value = 1422
$ws.Cells.Item(117, 6).Value = '''2019-03-12'
$ws.Cells.Item(118, 6).Value = '''2019-03-12'
$ws.Cells.Item(119, 6).Value = '''2019-03-12'
$ws.Cells.Item(120, 6).Value = '''2019-03-12'
$ws.Cells.Item(121, 6).Value = '''2019-03-12'
